$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4800
$ws.Range("I18").Value = 4800
$ws.Range("K18").Value = 4800
$ws.Range("M18").Value = -4516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 144.22223
$ws.Range("I33").Value = 130.57143
$ws.Range("K33").Value = 130.57143
$ws.Range("M33").Value = 98.42857000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4660
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3500
$ws.Range("I64").Value = 3500
$ws.Range("K64").Value = 3500
$ws.Range("M64").Value = -3252

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3500
$ws.Range("I67").Value = 3500
$ws.Range("K67").Value = 3500
$ws.Range("M67").Value = -2642

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2618.3667
$ws.Range("I137").Value = 904.3333
$ws.Range("J137").Value = 5189.4165
$ws.Range("K137").Value = 2712.9999
$ws.Range("L137").Value = 15568.2495
$ws.Range("M137").Value = -162.9998999999998
$ws.Range("N137").Value = -20668.2495

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6277.61
$ws.Range("I138").Value = 1260.3334
$ws.Range("J138").Value = 9172.191999999999
$ws.Range("K138").Value = 3781.0002
$ws.Range("L138").Value = 27516.576
$ws.Range("M138").Value = 1358.9998
$ws.Range("N138").Value = -37796.576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2599.75
$ws.Range("I45").Value = 2299.6667
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 2299.6667
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -1922.6667
$ws.Range("N45").Value = -4254

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 680.6429000000001
$ws.Range("J74").Value = 797
$ws.Range("L74").Value = 797
$ws.Range("N74").Value = -2545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 680.6429000000001
$ws.Range("J77").Value = 797
$ws.Range("L77").Value = 3985
$ws.Range("N77").Value = -12721

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 104
$ws.Range("I11").Value = 104
$ws.Range("K11").Value = 104
$ws.Range("M11").Value = 36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1142.7222
$ws.Range("I31").Value = 1050.3334
$ws.Range("J31").Value = 1327.5
$ws.Range("K31").Value = 1050.3334
$ws.Range("L31").Value = 1327.5
$ws.Range("M31").Value = -755.3334
$ws.Range("N31").Value = -1917.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1142.7222
$ws.Range("I34").Value = 1050.3334
$ws.Range("J34").Value = 1327.5
$ws.Range("K34").Value = 1050.3334
$ws.Range("L34").Value = 1327.5
$ws.Range("M34").Value = -848.3334
$ws.Range("N34").Value = -1731.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 38989.668
$ws.Range("J41").Value = 38984.5
$ws.Range("L41").Value = 38984.5
$ws.Range("N41").Value = -39840.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 50040
$ws.Range("J50").Value = 50040
$ws.Range("L50").Value = 50040
$ws.Range("N50").Value = -51290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4292.5
$ws.Range("I132").Value = 4102.8887
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 12308.6661
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -9778.666100000002
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 792
$ws.Range("I5").Value = 779
$ws.Range("J5").Value = 798.5
$ws.Range("K5").Value = 2337
$ws.Range("L5").Value = 2395.5
$ws.Range("M5").Value = -2225
$ws.Range("N5").Value = -2619.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 191995
$ws.Range("J37").Value = 191995
$ws.Range("L37").Value = 575985
$ws.Range("N37").Value = -576209

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 237.9
$ws.Range("I107").Value = 99
$ws.Range("J107").Value = 253.33333
$ws.Range("K107").Value = 297
$ws.Range("L107").Value = 759.99999
$ws.Range("M107").Value = 1623
$ws.Range("N107").Value = -4599.99999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3932.6667
$ws.Range("J132").Value = 3932.6667
$ws.Range("L132").Value = 35394.0003
$ws.Range("N132").Value = -40454.0003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 792
$ws.Range("I135").Value = 779
$ws.Range("J135").Value = 798.5
$ws.Range("K135").Value = 7011
$ws.Range("L135").Value = 7186.5
$ws.Range("M135").Value = -4476
$ws.Range("N135").Value = -12256.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 40000
$ws.Range("J20").Value = 40000
$ws.Range("L20").Value = 40000
$ws.Range("N20").Value = -40490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 19998
$ws.Range("I24").Value = 19990
$ws.Range("J24").Value = 20000
$ws.Range("K24").Value = 19990
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = -19817
$ws.Range("N24").Value = -20346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 34920.5
$ws.Range("J46").Value = 34920.5
$ws.Range("L46").Value = 34920.5
$ws.Range("N46").Value = -35232.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4896.6665
$ws.Range("I126").Value = 4896.6665
$ws.Range("K126").Value = 14689.9995
$ws.Range("M126").Value = -12219.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 864999.4399999999
$ws.Range("J43").Value = 864999.4399999999
$ws.Range("L43").Value = 864999.4399999999
$ws.Range("N43").Value = -865385.4399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1252.7273
$ws.Range("I82").Value = 1412.7142
$ws.Range("J82").Value = 972.75
$ws.Range("K82").Value = 1412.7142
$ws.Range("L82").Value = 972.75
$ws.Range("M82").Value = -1051.7142
$ws.Range("N82").Value = -1694.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1252.7273
$ws.Range("I85").Value = 1412.7142
$ws.Range("J85").Value = 972.75
$ws.Range("K85").Value = 1412.7142
$ws.Range("L85").Value = 972.75
$ws.Range("M85").Value = -164.7141999999999
$ws.Range("N85").Value = -3468.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1880
$ws.Range("J100").Value = 1950
$ws.Range("L100").Value = 1950
$ws.Range("N100").Value = -3032

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2131.5518
$ws.Range("I132").Value = 1545.4762
$ws.Range("K132").Value = 4636.4286
$ws.Range("M132").Value = -2106.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1094.2
$ws.Range("J126").Value = 809
$ws.Range("L126").Value = 2427
$ws.Range("M126").Value = -1383.0002
$ws.Range("N126").Value = -7367

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 767.94116
$ws.Range("I136").Value = 537.2
$ws.Range("J136").Value = 2498.5
$ws.Range("K136").Value = 1611.6
$ws.Range("L136").Value = 7495.5
$ws.Range("M136").Value = 938.3999999999999
$ws.Range("N136").Value = -12595.5
